$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy formatting from the last existing row (row 38) down to the new row 39
# so the date/time number formats match the rest of the column.
$null = $ws.Range("A38:D38").Copy()
$null = $ws.Range("A39:D39").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New entry for 2019-08-20 10:33 AM, 74 pages, word count 21139
$ws.Range("A39").Value = 43697
$ws.Range("B39").Value = 0.43958333333333338
$ws.Range("C39").Value = 74
$ws.Range("D39").Value = 21139

# Move the active selection to D40, matching where the user continued typing
$null = $ws.Range("D40").Select()
